$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix typo'd model names -- these labels repeat at the top of each of the
# four mini-tables (rows 2, 25, 29, 33) that feed the two charts.
foreach ($row in 2, 25, 29, 33) {
    $ws.Range("F$row").Value = "YOLOv8n"
    $ws.Range("H$row").Value = "YOLOv8 pt-mt"
}

# Restore the view state (scroll position / active selection) recorded in the saved file
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("O5").Select()
